$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.619.81'
$ws.Range('E2').Value = '  -6.38%  '
$ws.Range('D3').Value = '2.921.96'
$ws.Range('E3').Value = '  -8.20%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.90'
$ws.Range('E5').Value = '  -7.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.13'
$ws.Range('E6').Value = '  -10.27%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.544'
$ws.Range('E8').Value = '  -7.51%  '
$ws.Range('D9').Value = '2.923.86'
$ws.Range('E9').Value = '  -8.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.109'
$ws.Range('E10').Value = '  -7.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.07'
$ws.Range('E11').Value = '  -9.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  -8.85%  '
$ws.Range('D13').Value = '3.447.72'
$ws.Range('E13').Value = '  -8.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.123'
$ws.Range('E14').Value = '  -3.79%  '
$ws.Range('D15').Value = '60.884.35'
$ws.Range('E15').Value = '  -6.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.08'
$ws.Range('E16').Value = '  -9.53%  '
$ws.Range('D17').Value = '2.940.44'
$ws.Range('E17').Value = '  -7.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000142'
$ws.Range('E18').Value = '  -9.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.01'
$ws.Range('E19').Value = '  -5.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '375.20'
$ws.Range('E20').Value = '  -8.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.63'
$ws.Range('E21').Value = '  -9.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.43'
$ws.Range('E22').Value = '  -10.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.94'
$ws.Range('E24').Value = '  -7.77%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.088.28'
$ws.Range('E25').Value = '  -7.63%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.456'
$ws.Range('E26').Value = '  -7.12%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.183'
$ws.Range('E27').Value = '  -8.89%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.14'
$ws.Range('E29').Value = '  -8.19%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0893'
$ws.Range('E30').Value = '  -14.33%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.67'
$ws.Range('E32').Value = '  -8.73%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.71'
$ws.Range('E33').Value = '  -7.87%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '156.35'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.85'
$ws.Range('E35').Value = '  -8.19%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.46'
$ws.Range('E36').Value = '  -10.07%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.04'
$ws.Range('E37').Value = '  -8.99%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.24'
$ws.Range('E38').Value = '  -8.76%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.50'
$ws.Range('E39').Value = '  -11.72%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.384.85'
$ws.Range('E40').Value = '  -12.42%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.41'
$ws.Range('E41').Value = '  -6.56%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.79'
$ws.Range('E42').Value = '  -7.93%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.59'
$ws.Range('E43').Value = '  -10.99%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.649'
$ws.Range('E44').Value = '  -8.58%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0581'
$ws.Range('E45').Value = '  -8.23%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0239'
$ws.Range('E47').Value = '  -8.65%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.48'
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.82'
$ws.Range('E49').Value = '  -12.69%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0932'
$ws.Range('E50').Value = '  -5.60%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.20'
$ws.Range('E51').Value = '  -10.82%  '
